$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Question 1 block (rows 6-11): COGS now absorbs the advertising expense,
# margins recompute against RFC Gross Sales Revenue (I4) instead of Net (I6)
# ---------------------------------------------------------------------------
$ws.Range("D6").Formula = "=I4-I8"
$ws.Range("D7").Formula = "=D6/I4 * 100"
$ws.Range("D8").Formula = "=I9+I10+I13+I5"
$ws.Range("D9").Formula = "=I4-I8-D8"
$ws.Range("D10").Formula = "=D9/I4 * 100"
$ws.Range("D11").Formula = "=D9/D8 *100"

# ---------------------------------------------------------------------------
# Question 2 block (rows 22-27): COGS increases with volume, trade promotion
# now feeds MSE, margins recompute against I20 instead of I22
# ---------------------------------------------------------------------------
$ws.Range("D24").Formula = "=I25+I26+I29+I21"
$ws.Range("I24").Formula = "=I8 + I8*(8/15)"
$ws.Range("J24").Value = "changes because more goods therefore more COGs"
$ws.Range("D25").Formula = "=I20-I24-D24"
$ws.Range("D26").Formula = "=(D25/I20)*100"

# ---------------------------------------------------------------------------
# Question 3 block (rows 36-45): COGS now absorbs trade promotion,
# margins recompute against I36 instead of I38
# ---------------------------------------------------------------------------
$ws.Range("D38").Formula = "=I36-I40"
$ws.Range("D39").Formula = "=(D38/I36)*100"
$ws.Range("D40").Formula = "=I41+I42+I45+I37"
$ws.Range("D41").Formula = "=I36-I40-D40"
$ws.Range("D42").Formula = "=(D41/I36)*100"
$ws.Range("E42").Value = "%"
$ws.Range("D43").Formula = "=(D41/D40)*100"
$ws.Range("E43").Value = "%"

# Clear the old boxed note (C45:E45) and replace with an unstyled label
$ws.Range("C45:E45").Style = "Normal"
$ws.Range("D45").ClearContents()
$ws.Range("E45").ClearContents()
$ws.Range("C45").Value = "New Req Sales Revenue = (NMC + MSE + COGS) / (1 -Increased Trade Percentage)"

# Give D48 the currency-ish numeric style (s="5") used elsewhere, left blank
$ws.Range("D48").Formula = ""
$ws.Range("D48").Style = $ws.Range("D41").Style

# ---------------------------------------------------------------------------
# Replace old row 49 with the new Question-2 COGS scratch-work (rows 50-51)
# and the relocated "Increase in $ sales needed" row (row 54)
# ---------------------------------------------------------------------------
$ws.Range("C49").ClearContents()
$ws.Range("D49").ClearContents()
$ws.Range("E49").ClearContents()

$ws.Range("C50").Value = "600,000 + 20*0.15*units sold"
$ws.Range("C51").Value = "900,000 = 9*goods sold - 600,000"
$ws.Range("D51").Value = "166,666 need to be sold"

$ws.Range("C54").Value = "Increase in $ sales needed"
$ws.Range("D54").Value = 11.11
$ws.Range("E54").Value = "%"
